# Libro.xlsx edit: update id value, remove the blank gap row, and move
# the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a blank row 3 (only row 4 carried a lone styled, empty
# cell E4). Deleting the empty row 3 shifts row 4 up to row 3 (E4 -> E3),
# shrinking the used range from A1:E4 to A1:E3.
$ws.Rows(3).Delete() | Out-Null

# A2 ("id" column) changes from 12345 to 1.
$ws.Range("A2").Value = 1

# Move the active selection to D14.
$ws.Range("D14").Select() | Out-Null
